$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row (row 1) labels to the new taxonomy.
# Columns A-F and O keep their original labels; G,H,I,J,K,L,M,N,P are renamed.
# (Order matches the new shared-string table layout.)
$ws.Range("G1").Value = "articlePages"
$ws.Range("H1").Value = "articleSeq"
$ws.Range("M1").Value = "fileName1"
$ws.Range("I1").Value = "authorGivenname1"
$ws.Range("J1").Value = "authorFamilyname1"
$ws.Range("K1").Value = "authorGivenname2"
$ws.Range("L1").Value = "authorFamilyname2"
$ws.Range("N1").Value = "galleyLabel1"
$ws.Range("P1").Value = "galleyLocale1"

# Move the active selection to I9, matching the saved view state.
$ws.Range("I9").Select()
